# Add a new student record (Luke) as row 5 on the "Class3" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class3")

# Every other populated cell on this sheet stores its numeric-looking values
# as literal text, so pre-format the numeric-looking text columns (Roll.No,
# Date of Birth, Class, Mark1-Mark5) as Text before writing to them --
# otherwise Excel would silently reinterpret strings like "21" or
# "12-01-2000" as a real number / date instead of keeping the literal text.
$ws.Range("B5:D5").NumberFormat = "@"
$ws.Range("F5:J5").NumberFormat = "@"

$ws.Range("A5").Value = "Luke"
$ws.Range("B5").Value = "21"
$ws.Range("C5").Value = "12-01-2000"
$ws.Range("D5").Value = "3"
$ws.Range("E5").Value = "USA"
$ws.Range("F5").Value = "77"
$ws.Range("G5").Value = "77"
$ws.Range("H5").Value = "100"
$ws.Range("I5").Value = "77"
$ws.Range("J5").Value = "77"
$ws.Range("K5").Value = 408
$ws.Range("L5").Value = "B"
